$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("REVIEW-SHEET")
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

$ws2.Range("A3").Value = "v1.1"

$ws1.Range("H2").Value = "not applicable"
$ws1.Range("H3").Value = "not applicable"
$ws1.Range("H4").Value = "not applicable"

$ws2.Range("B3").Value = "Omar Sherif"
$ws2.Range("C3").Value = "All comments closed as it is not applicaple "
$ws2.Range("D3").Value = 45772
